# Apply ETL-consolidator refactor: replace data rows 2-11 with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(86449, "Sra. Helena Silva",       "Marketing",         "Viagem de negócios",   1, 45081,  7280.22),
    @(87110, "Helena da Mota",          "Recursos Humanos",  "Consulta médica",      4, 45086,  2828.96),
    @(69674, "Isadora Ferreira",        "Jurídico",          "Viagem de negócios",   2, 45091,  3330.47),
    @(23911, "Juliana Jesus",           "Jurídico",          "Doença",               1, 45098,  8593.27),
    @(49118, "Igor Costa",              "Jurídico",          "Doença",               3, 45078,  4932.81),
    @(3555,  "Nathan da Mata",          "Recursos Humanos",  "Consulta médica",      3, 45102,  8658.65),
    @(5708,  "Lucas Gabriel da Luz",    "Recursos Humanos",  "Problemas pessoais",   5, 45098,  8199.1),
    @(19364, "Dra. Olivia Costela",     "Marketing",         "Outros",               6, 45097,  3030),
    @(87625, "Ana Beatriz Freitas",     "Recursos Humanos",  "Outros",               6, 45088,  3041.28),
    @(68104, "Luiz Gustavo Nascimento", "Jurídico",          "Outros",               8, 45085,  3653)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
